$wb = $excel.ActiveWorkbook
$wsTraining = $wb.Worksheets.Item("Training Dashboard")
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# Training Dashboard: refresh "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I)
# for every data row (3-27). The report was re-run 8 days later (16-Sep-2025
# instead of 08-Sep-2025), so every countdown value drops by 8.
# ---------------------------------------------------------------------------
for ($row = 3; $row -le 27; $row++) {
    $hCell = $wsTraining.Cells.Item($row, 8)
    $hCell.Value2 = $hCell.Value2 - 8
    $wsTraining.Cells.Item($row, 9).Value2 = "'16-Sep-2025"
}

# Row 18 (LOTO (SOPs)) now falls within the "about to expire" window and
# flips from VALID to NOT VALID, picking up the same pink highlight used by
# the other NOT VALID rows (19-22).
$wsTraining.Cells.Item(18, 10).Value2 = "NOT VALID"
$wsTraining.Range("A18:K18").Interior.Color = 13551615
$wsTraining.Range("A18:K18").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Header styling: bold white text on the dark-blue fill for the title row and
# the column-header row (both sheets), matching the refreshed font table.
# ---------------------------------------------------------------------------
$wsTraining.Range("A1").Font.Bold = $true
$wsTraining.Range("A1").Font.Size = 11
$wsTraining.Range("A1").Font.Color = 16777215
$wsTraining.Range("A2:K2").Font.Bold = $true
$wsTraining.Range("A2:K2").Font.Color = 16777215

$wsExam.Range("A1").Font.Bold = $true
$wsExam.Range("A1").Font.Size = 11
$wsExam.Range("A1").Font.Color = 16777215
$wsExam.Range("A2:G2").Font.Bold = $true
$wsExam.Range("A2:G2").Font.Color = 16777215

# ---------------------------------------------------------------------------
# Exam Dashboard: reword the per-exam comment and widen the comments column.
# ---------------------------------------------------------------------------
for ($row = 3; $row -le 12; $row++) {
    $wsExam.Cells.Item($row, 5).Value2 = "date is valid"
}

# (COM's ColumnWidth setter adds the usual 5/6-character grid padding on
# save, so back it off here to land on a stored width of exactly 15.)
$wsExam.Columns.Item(5).ColumnWidth = 14.166666666666666
